$d = $word.ActiveDocument

# The document currently has a single empty paragraph (its paragraph mark
# is already formatted en-US). Fill it with "Test 1 desktop", then add a
# new paragraph after it with "Test 2 Laptop" -- matching the original,
# reverted content.
$para1 = $d.Paragraphs(1)
$para1.Range.Text = "Test 1 desktop"
$para1.Range.LanguageID = "en-US"

$end = $d.Content
$end.Collapse(0)
$end.InsertParagraphAfter()
$end.Collapse(0)
$end.Text = "Test 2 Laptop"
$end.LanguageID = "en-US"

